$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 10.23061133333333
$ws.Cells.Item(2, 8).Value = 30.691834
$ws.Cells.Item(2, 9).Value = 0.4855635428718841
$ws.Cells.Item(2, 10).Value = 0.4855635428718841
$ws.Cells.Item(2, 13).Value = 9.546140333333334
$ws.Cells.Item(2, 14).Value = 28.638421
$ws.Cells.Item(2, 15).Value = 0.587227294878132
$ws.Cells.Item(2, 16).Value = 0.587227294878132
$ws.Cells.Item(2, 17).Value = 97.66285148379045
$ws.Cells.Item(2, 18).Value = 878.965663354114
$ws.Cells.Item(2, 19).Value = 0.2851361657720983
$ws.Cells.Item(2, 20).Value = 0.2851361657720983

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 10.23061133333333
$ws.Cells.Item(3, 8).Value = 30.691834
$ws.Cells.Item(3, 9).Value = 0.4855635428718841
$ws.Cells.Item(3, 10).Value = 0.4855635428718841
$ws.Cells.Item(3, 15).Value = 0.2496684258894083
$ws.Cells.Item(3, 16).Value = 0.2496684258894083
$ws.Cells.Item(3, 17).Value = 41.52281511861489
$ws.Cells.Item(3, 18).Value = 373.705336067534
$ws.Cells.Item(3, 19).Value = 0.1212298854181075
$ws.Cells.Item(3, 20).Value = 0.1212298854181075

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 10.23061133333333
$ws.Cells.Item(4, 8).Value = 30.691834
$ws.Cells.Item(4, 9).Value = 0.4855635428718841
$ws.Cells.Item(4, 10).Value = 0.4855635428718841
$ws.Cells.Item(4, 13).Value = 2.210442
$ws.Cells.Item(4, 14).Value = 6.631326
$ws.Cells.Item(4, 15).Value = 0.1359745227725727
$ws.Cells.Item(4, 16).Value = 0.1359745227725727
$ws.Cells.Item(4, 17).Value = 22.614172976876
$ws.Cells.Item(4, 18).Value = 203.527556791884
$ws.Cells.Item(4, 19).Value = 0.06602427101776406
$ws.Cells.Item(4, 20).Value = 0.06602427101776406

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 10.23061133333333
$ws.Cells.Item(5, 8).Value = 30.691834
$ws.Cells.Item(5, 9).Value = 0.4855635428718841
$ws.Cells.Item(5, 10).Value = 0.4855635428718841
$ws.Cells.Item(5, 13).Value = 0.4410293333333333
$ws.Cells.Item(5, 14).Value = 1.323088
$ws.Cells.Item(5, 15).Value = 0.02712975645988715
$ws.Cells.Item(5, 16).Value = 0.02712975645988715
$ws.Cells.Item(5, 17).Value = 4.511999695932444
$ws.Cells.Item(5, 18).Value = 40.607997263392
$ws.Cells.Item(5, 19).Value = 0.01317322066391419
$ws.Cells.Item(5, 20).Value = 0.01317322066391419

# Row 6
$ws.Cells.Item(6, 9).Value = 0.4164864079521221
$ws.Cells.Item(6, 10).Value = 0.4164864079521222
$ws.Cells.Item(6, 13).Value = 9.546140333333334
$ws.Cells.Item(6, 14).Value = 28.638421
$ws.Cells.Item(6, 15).Value = 0.587227294878132
$ws.Cells.Item(6, 16).Value = 0.587227294878132
$ws.Cells.Item(6, 17).Value = 83.76916018914876
$ws.Cells.Item(6, 18).Value = 753.922441702339
$ws.Cells.Item(6, 19).Value = 0.2445721866952348
$ws.Cells.Item(6, 20).Value = 0.2445721866952348

# Row 7
$ws.Cells.Item(7, 9).Value = 0.4164864079521221
$ws.Cells.Item(7, 10).Value = 0.4164864079521222
$ws.Cells.Item(7, 15).Value = 0.2496684258894083
$ws.Cells.Item(7, 16).Value = 0.2496684258894083
$ws.Cells.Item(7, 19).Value = 0.1039835058777403
$ws.Cells.Item(7, 20).Value = 0.1039835058777403

# Row 8
$ws.Cells.Item(8, 9).Value = 0.4164864079521221
$ws.Cells.Item(8, 10).Value = 0.4164864079521222
$ws.Cells.Item(8, 13).Value = 2.210442
$ws.Cells.Item(8, 14).Value = 6.631326
$ws.Cells.Item(8, 15).Value = 0.1359745227725727
$ws.Cells.Item(8, 16).Value = 0.1359745227725727
$ws.Cells.Item(8, 17).Value = 19.397040429026
$ws.Cells.Item(8, 18).Value = 174.573363861234
$ws.Cells.Item(8, 19).Value = 0.05663154056255281
$ws.Cells.Item(8, 20).Value = 0.05663154056255282

# Row 9
$ws.Cells.Item(9, 9).Value = 0.4164864079521221
$ws.Cells.Item(9, 10).Value = 0.4164864079521222
$ws.Cells.Item(9, 13).Value = 0.4410293333333333
$ws.Cells.Item(9, 14).Value = 1.323088
$ws.Cells.Item(9, 15).Value = 0.02712975645988715
$ws.Cells.Item(9, 16).Value = 0.02712975645988715
$ws.Cells.Item(9, 17).Value = 3.870114578465777
$ws.Cells.Item(9, 18).Value = 34.831031206192
$ws.Cells.Item(9, 19).Value = 0.01129917481659428
$ws.Cells.Item(9, 20).Value = 0.01129917481659428

# Row 10
$ws.Cells.Item(10, 7).Value = 2.034752
$ws.Cells.Item(10, 8).Value = 6.104255999999999
$ws.Cells.Item(10, 9).Value = 0.09657305490303886
$ws.Cells.Item(10, 10).Value = 0.09657305490303887
$ws.Cells.Item(10, 13).Value = 9.546140333333334
$ws.Cells.Item(10, 14).Value = 28.638421
$ws.Cells.Item(10, 15).Value = 0.587227294878132
$ws.Cells.Item(10, 16).Value = 0.587227294878132
$ws.Cells.Item(10, 17).Value = 19.42402813553067
$ws.Cells.Item(10, 18).Value = 174.816253219776
$ws.Cells.Item(10, 19).Value = 0.05671033378882883
$ws.Cells.Item(10, 20).Value = 0.05671033378882884

# Row 11
$ws.Cells.Item(11, 7).Value = 2.034752
$ws.Cells.Item(11, 8).Value = 6.104255999999999
$ws.Cells.Item(11, 9).Value = 0.09657305490303886
$ws.Cells.Item(11, 10).Value = 0.09657305490303887
$ws.Cells.Item(11, 15).Value = 0.2496684258894083
$ws.Cells.Item(11, 16).Value = 0.2496684258894083
$ws.Cells.Item(11, 17).Value = 8.258414708117332
$ws.Cells.Item(11, 18).Value = 74.325732373056
$ws.Cells.Item(11, 19).Value = 0.02411124260097311
$ws.Cells.Item(11, 20).Value = 0.02411124260097312

# Row 12
$ws.Cells.Item(12, 7).Value = 2.034752
$ws.Cells.Item(12, 8).Value = 6.104255999999999
$ws.Cells.Item(12, 9).Value = 0.09657305490303886
$ws.Cells.Item(12, 10).Value = 0.09657305490303887
$ws.Cells.Item(12, 13).Value = 2.210442
$ws.Cells.Item(12, 14).Value = 6.631326
$ws.Cells.Item(12, 15).Value = 0.1359745227725727
$ws.Cells.Item(12, 16).Value = 0.1359745227725727
$ws.Cells.Item(12, 17).Value = 4.497701280384
$ws.Cells.Item(12, 18).Value = 40.479311523456
$ws.Cells.Item(12, 19).Value = 0.01313147505313017
$ws.Cells.Item(12, 20).Value = 0.01313147505313017

# Row 13
$ws.Cells.Item(13, 7).Value = 2.034752
$ws.Cells.Item(13, 8).Value = 6.104255999999999
$ws.Cells.Item(13, 9).Value = 0.09657305490303886
$ws.Cells.Item(13, 10).Value = 0.09657305490303887
$ws.Cells.Item(13, 13).Value = 0.4410293333333333
$ws.Cells.Item(13, 14).Value = 1.323088
$ws.Cells.Item(13, 15).Value = 0.02712975645988715
$ws.Cells.Item(13, 16).Value = 0.02712975645988715
$ws.Cells.Item(13, 17).Value = 0.8973853180586665
$ws.Cells.Item(13, 18).Value = 8.076467862528
$ws.Cells.Item(13, 19).Value = 0.002620003460106754
$ws.Cells.Item(13, 20).Value = 0.002620003460106755

# Row 14
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.02901266666666667
$ws.Cells.Item(14, 8).Value = 0.087038
$ws.Cells.Item(14, 9).Value = 0.001376994272954919
$ws.Cells.Item(14, 10).Value = 0.001376994272954919
$ws.Cells.Item(14, 13).Value = 9.546140333333334
$ws.Cells.Item(14, 14).Value = 28.638421
$ws.Cells.Item(14, 15).Value = 0.587227294878132
$ws.Cells.Item(14, 16).Value = 0.587227294878132
$ws.Cells.Item(14, 17).Value = 0.2769589874442223
$ws.Cells.Item(14, 18).Value = 2.492630886998
$ws.Cells.Item(14, 19).Value = 0.0008086086219699969
$ws.Cells.Item(14, 20).Value = 0.0008086086219699969

# Row 15
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.02901266666666667
$ws.Cells.Item(15, 8).Value = 0.087038
$ws.Cells.Item(15, 9).Value = 0.001376994272954919
$ws.Cells.Item(15, 10).Value = 0.001376994272954919
$ws.Cells.Item(15, 15).Value = 0.2496684258894083
$ws.Cells.Item(15, 16).Value = 0.2496684258894083
$ws.Cells.Item(15, 17).Value = 0.1177532363264445
$ws.Cells.Item(15, 18).Value = 1.059779126938
$ws.Cells.Item(15, 19).Value = 0.0003437919925873847
$ws.Cells.Item(15, 20).Value = 0.0003437919925873847

# Row 16
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.02901266666666667
$ws.Cells.Item(16, 8).Value = 0.087038
$ws.Cells.Item(16, 9).Value = 0.001376994272954919
$ws.Cells.Item(16, 10).Value = 0.001376994272954919
$ws.Cells.Item(16, 13).Value = 2.210442
$ws.Cells.Item(16, 14).Value = 6.631326
$ws.Cells.Item(16, 15).Value = 0.1359745227725727
$ws.Cells.Item(16, 16).Value = 0.1359745227725727
$ws.Cells.Item(16, 17).Value = 0.06413081693200001
$ws.Cells.Item(16, 18).Value = 0.5771773523880001
$ws.Cells.Item(16, 19).Value = 0.0001872361391256107
$ws.Cells.Item(16, 20).Value = 0.0001872361391256107

# Row 17
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.02901266666666667
$ws.Cells.Item(17, 8).Value = 0.087038
$ws.Cells.Item(17, 9).Value = 0.001376994272954919
$ws.Cells.Item(17, 10).Value = 0.001376994272954919
$ws.Cells.Item(17, 13).Value = 0.4410293333333333
$ws.Cells.Item(17, 14).Value = 1.323088
$ws.Cells.Item(17, 15).Value = 0.02712975645988715
$ws.Cells.Item(17, 17).Value = 0.01279543703822222
$ws.Cells.Item(17, 18).Value = 0.115158933344
$ws.Cells.Item(17, 19).Value = [double]"3.735751927192631E-05"
$ws.Cells.Item(17, 20).Value = [double]"3.735751927192631E-05"
